$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The burn-down actuals dropped sharply at row 16 (C16): instead of holding
# steady at 101, the remaining work fell to 11 and stayed there through the
# end of the sprint (C17, C18 follow the existing "carry forward" formula).
$ws.Range("C16").Formula = "=C15-90"

# Reflect where the author's cursor ended up after making the edit.
[void]$ws.Range("D16").Select()
